$wb = $excel.ActiveWorkbook

$wsC  = $wb.Worksheets.Item("C")
$wsC0 = $wb.Worksheets.Item("C0")

# Sheet "C": column B formulas change from pulling R_input!$C to
# referencing column A of the same sheet ('C'!$A.. & "init")
for ($r = 1; $r -le 60; $r++) {
    $wsC.Range("B$r").Formula = "=IF('C'!`$A$r=`"`",`"`",'C'!`$A$r&`"init`")"
}

# Sheet "C0": column B formulas change from referencing sheet C's column B
# back to pulling directly from R_input!$C
for ($r = 1; $r -le 60; $r++) {
    $wsC0.Range("B$r").Formula = "=IF(R_input!`$C$r=0,`"`",R_input!`$C$r)"
}

# Update the selections to match the recorded cursor positions
$wsC.Range("D6").Select()
$wsC0.Range("F7").Select()

$wb.Application.Calculate()

# Restore the originally active sheet (R_input) so the workbook-level
# active tab is unaffected by the selection changes above
$wb.Worksheets.Item("R_input").Activate()
